# Update "roberta + cnn" and "bert + bilstm" result:
# Restructure the "baseline" sheet into a consolidated table with a running
# row-index column, absorb the "roberta + lstm" row that used to be
# duplicated on the "resampling" sheet, and append four new model rows
# (deberta+lstm, bert+cnn, robert+cnn, bert+bilstm) with two spacer rows,
# matching the new layout recorded in the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "baseline"
$ws2 = $wb.Worksheets.Item(2)   # "resampling"

# ---------------------------------------------------------------------------
# 1) "baseline" sheet: insert a new leading column that will hold a running
#    row index (1,2,3,...), shifting the existing B:I data to C:J.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).Insert()

$ws1.Columns.Item(1).ColumnWidth = 3.2857142857142856   # -> stored width 4
$ws1.Columns.Item(2).ColumnWidth = 12.426339285714286   # -> stored width ~13.14

# Row index numbers for the existing rows 2-5 (bert, roberta, deberta, bert+lstm)
$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(5, 1).Value = 4

# ---------------------------------------------------------------------------
# 2) Append the "roberta + lstm" row (previously duplicated on "resampling"
#    row 6) as row 6 on "baseline", then four new model rows, then spacers.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 6;  Idx = 5;  Name = "roberta + lstm"; Vals = @(1.4501999999999999, 0.65380000000000005, 0.64070000000000005, 0.65380000000000005, 0.63419999999999999, 0.30430000000000001, 0.31640000000000001, 0.54590000000000005) },
    @{ Row = 7;  Idx = 6;  Name = "deberta + lstm";  Vals = @(2.4902000000000002, 0.46920000000000001, 0.22650000000000001, 0.46920000000000001, 0.30409999999999998, 0.0152, 0.0218, 0.0284) },
    @{ Row = 8;  Idx = 7;  Name = "bert + cnn";      Vals = @(1.5245, 0.67849999999999999, 0.68240000000000001, 0.67849999999999999, 0.66649999999999998, 0.41959999999999997, 0.42309999999999998, 0.57599999999999996) },
    @{ Row = 9;  Idx = 8;  Name = "robert + cnn";    Vals = @(1.4946999999999999, 0.68100000000000005, 0.68710000000000004, 0.68100000000000005, 0.66990000000000005, 0.40720000000000001, 0.40849999999999997, 0.58109999999999995) },
    @{ Row = 11; Idx = 10; Name = "bert + bilstm";   Vals = @(1.5226, 0.68220000000000003, 0.6845, 0.68220000000000003, 0.67200000000000004, 0.40410000000000001, 0.41389999999999999, 0.58230000000000004) }
)

foreach ($r in $newRows) {
    $ws1.Cells.Item($r.Row, 1).Value = $r.Idx
    $ws1.Cells.Item($r.Row, 2).Value = $r.Name
    for ($c = 0; $c -lt $r.Vals.Length; $c++) {
        $cell = $ws1.Cells.Item($r.Row, 3 + $c)
        $cell.Value = $r.Vals[$c]
        $cell.NumberFormat = "0.0000"
    }
}

# Spacer rows: only the running index in column A, nothing else.
$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(12, 1).Value = 11
$ws1.Cells.Item(13, 1).Value = 12

# ---------------------------------------------------------------------------
# 3) "resampling" sheet: remove the now-redundant "roberta + lstm" row (row 6)
#    and give the sheet a page setup (paper size 9 / portrait).
# ---------------------------------------------------------------------------
$ws2.Rows.Item(6).Delete()

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("C15").Select()

# ---------------------------------------------------------------------------
# 4) Make "baseline" the active sheet/tab (it was "resampling" before).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("I14").Select()
